$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from the "8.2. Ferramentas..." heading to the
#    end of the "Google Chrome" browser/version list item (right after the
#    version text, before the paragraph mark). We do this first (while
#    paragraph positions are still the "before deletion" layout) by inserting
#    a temporary placeholder character at the target insertion point, wrapping
#    a bookmark named "_GoBack" around it (re-adding a bookmark with this
#    reserved name automatically removes/replaces any existing one elsewhere
#    in the document, mirroring Word's single-instance "_GoBack" behaviour),
#    and then deleting the placeholder again so the bookmark collapses to the
#    correct zero-length insertion point.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*141.0.7390.123*") {
        $r = $p.Range.Duplicate
        $r.MoveEnd(1, -1) | Out-Null
        $r.Collapse(0)
        $r.InsertAfter("X")
        $d.Bookmarks.Add("_GoBack", $r) | Out-Null
        $r.Text = ""
        break
    }
}

# 2) Delete the whole "Mozilla Firefox: Versão 121.0 ou superior" list item
#    paragraph (including its paragraph mark).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Mozilla Firefox*") {
        $p.Range.Delete()
        break
    }
}
